# Header and Flip Sign Dimension and Location KPIs running

$wb = $excel.ActiveWorkbook

# --- Fix a typo on the "Header Positions" sheet: "Cigarette Positions" -> "Cigarettes Positions"
$wsHeader = $wb.Worksheets.Item("Header Positions")
$wsHeader.Range("B1").Value = "Cigarettes Positions"

# --- Rename the last sheet from "POS smokeless" to "Flip Sign Positions"
$wsFlip = $wb.Worksheets.Item("POS smokeless")
$wsFlip.Name = "Flip Sign Positions"

# --- Move the selection on "Header Positions" from C13 to D21 (no longer the active tab)
[void]$wsHeader.Range("D21").Select()

# --- Make "Flip Sign Positions" the active/selected sheet and tab
[void]$wsFlip.Activate()
